$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 17.62433929335113
$ws.Range("C2").Value = 10.93012825351122
$ws.Range("E2").Value = 11.65254030341503
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 3.635237828742222
$ws.Range("L2").Value = 9.771818357878768
$ws.Range("O2").Value = 21.783745325392
$ws.Range("B3").Value = 16.97238938020411
$ws.Range("C3").Value = 10.70262724041649
$ws.Range("E3").Value = 11.71219608399484
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 3.637503369411146
$ws.Range("L3").Value = 9.73886407234202
$ws.Range("O3").Value = 21.94809718031337
$ws.Range("B4").Value = 16.56012919542524
$ws.Range("C4").Value = 10.56044015581627
$ws.Range("E4").Value = 11.75123777084512
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 3.638966683031612
$ws.Range("L4").Value = 9.720197387447698
$ws.Range("O4").Value = 22.05634916396239
$ws.Range("B5").Value = 16.38935770930619
$ws.Range("C5").Value = 10.501932405971
$ws.Range("E5").Value = 11.7677543857409
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 3.639581225781412
$ws.Range("L5").Value = 9.712989811305377
$ws.Range("O5").Value = 22.10230152703428
$ws.Range("B6").Value = 16.36084117656275
$ws.Range("C6").Value = 10.49218499937948
$ws.Range("E6").Value = 11.77053361014463
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 3.639684372969924
$ws.Range("L6").Value = 9.711817254299502
$ws.Range("O6").Value = 22.11004274790297
$ws.Range("B7").Value = 16.55783700619521
$ws.Range("C7").Value = 10.55965330314231
$ws.Range("E7").Value = 11.75145806259686
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 3.638974897075788
$ws.Range("L7").Value = 9.720098560538947
$ws.Range("O7").Value = 22.05696145716558
$ws.Range("B8").Value = 17.40216287468084
$ws.Range("C8").Value = 10.85224207202342
$ws.Range("E8").Value = 11.67260889373546
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 3.636004025013389
$ws.Range("L8").Value = 9.760133596809158
$ws.Range("O8").Value = 21.83888649128818
$ws.Range("B9").Value = 18.95388323225754
$ws.Range("C9").Value = 11.40352893379448
$ws.Range("E9").Value = 11.53712617022188
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.630748798001679
$ws.Range("L9").Value = 9.850837334963121
$ws.Range("O9").Value = 21.46978281317794
$ws.Range("B10").Value = 20.01997720969961
$ws.Range("C10").Value = 11.79154642103919
$ws.Range("E10").Value = 11.44924972112825
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.627231836297645
$ws.Range("L10").Value = 9.924565716020464
$ws.Range("O10").Value = 21.23472029092688
$ws.Range("B11").Value = 20.48711007994854
$ws.Range("C11").Value = 11.96373794246101
$ws.Range("E11").Value = 11.41180426319936
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.625705771616263
$ws.Range("L11").Value = 9.959564825872855
$ws.Range("O11").Value = 21.13572366844301
$ws.Range("B12").Value = 20.66131109924978
$ws.Range("C12").Value = 12.02827300072829
$ws.Range("E12").Value = 11.39798853427294
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.625138443815024
$ws.Range("L12").Value = 9.973020606923983
$ws.Range("O12").Value = 21.09938569976714
$ws.Range("B13").Value = 20.62391536374
$ws.Range("C13").Value = 12.01440476814087
$ws.Range("E13").Value = 11.40094780308737
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.625260159235844
$ws.Range("L13").Value = 9.97011378583859
$ws.Range("O13").Value = 21.10716044450774
$ws.Range("B14").Value = 20.50149632209651
$ws.Range("C14").Value = 11.96906097907927
$ws.Range("E14").Value = 11.41066033717965
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.625658885898495
$ws.Range("L14").Value = 9.960667829269291
$ws.Range("O14").Value = 21.13271102412202
$ws.Range("B15").Value = 20.42615706316938
$ws.Range("C15").Value = 11.94119793415668
$ws.Range("E15").Value = 11.41665696305264
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.625904491052514
$ws.Range("L15").Value = 9.954908035330902
$ws.Range("O15").Value = 21.14851151235082
$ws.Range("B16").Value = 19.98907877923526
$ws.Range("C16").Value = 11.78020211145394
$ws.Range("E16").Value = 11.45174779979141
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.627333048903473
$ws.Range("L16").Value = 9.922307178278999
$ws.Range("O16").Value = 21.24135037371804
$ws.Range("B17").Value = 19.71628070274372
$ws.Range("C17").Value = 11.68029406985275
$ws.Range("E17").Value = 11.47392304478358
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.628228290196333
$ws.Range("L17").Value = 9.902676213712329
$ws.Range("O17").Value = 21.30034250587637
$ws.Range("B18").Value = 19.55770183742405
$ws.Range("C18").Value = 11.62242573216842
$ws.Range("E18").Value = 11.48691577622309
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.628750160769542
$ws.Range("L18").Value = 9.891522866005557
$ws.Range("O18").Value = 21.33501950366519
$ws.Range("B19").Value = 19.50372663304234
$ws.Range("C19").Value = 11.60276469175823
$ws.Range("E19").Value = 11.49135578170938
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.628928052788723
$ws.Range("L19").Value = 9.887770436757533
$ws.Range("O19").Value = 21.3468884458883
$ws.Range("B20").Value = 19.74549460778901
$ws.Range("C20").Value = 11.69097161695312
$ws.Range("E20").Value = 11.47153780619535
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.62813227118372
$ws.Range("L20").Value = 9.904751748580457
$ws.Range("O20").Value = 21.29398539112385
$ws.Range("B21").Value = 20.53752775104199
$ws.Range("C21").Value = 11.98239810923049
$ws.Range("E21").Value = 11.40779764807034
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.625541484061732
$ws.Range("L21").Value = 9.963436904635484
$ws.Range("O21").Value = 21.12517491914011
$ws.Range("B22").Value = 21.03942856953235
$ws.Range("C22").Value = 12.16893587055629
$ws.Range("E22").Value = 11.3682618872786
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.623909779956502
$ws.Range("L22").Value = 10.00296696767532
$ws.Range("O22").Value = 21.02155434444326
$ws.Range("B23").Value = 20.77303235578109
$ws.Range("C23").Value = 12.06975168777334
$ws.Range("E23").Value = 11.38916860390491
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.624775039772553
$ws.Range("L23").Value = 9.981763992507231
$ws.Range("O23").Value = 21.07624190914741
$ws.Range("B24").Value = 19.73229242249473
$ws.Range("C24").Value = 11.68614563314582
$ws.Range("E24").Value = 11.47261541288248
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.628175658995516
$ws.Range("L24").Value = 9.903812985185864
$ws.Range("O24").Value = 21.29685707085684
$ws.Range("B25").Value = 18.54641094534368
$ws.Range("C25").Value = 11.25715228370133
$ws.Range("E25").Value = 11.57172944188004
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.632109781660255
$ws.Range("L25").Value = 9.825028280552701
$ws.Range("O25").Value = 21.56332438713492
